$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (I1, J1), matching the formatting used by the other
# header cells (bold, bordered, centered) by copying H1's format.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data columns I and J, rows 2-5 (plain, unstyled numeric cells)
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 3

$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 6

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 2
